$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new EUR->ARS rate observation as row 77 (A:C), keeping every
# value as plain text so Excel doesn't auto-coerce the date/time strings
# into date/time serial numbers.
$ws.Range("A77:C77").NumberFormat = "@"

$ws.Range("A77").Value = "2025-10-14"
$ws.Range("B77").Value = "15:23:50"
$ws.Range("C77").Value = "1.00 EUR = 1,747.2536"
